# Scheduled-runner update: refresh market-board price snapshots (and
# derived profit columns) for a handful of leves across the Sheets
# workbook. Mirrors the "Sheets via scheduled runner" price-sync job:
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# (H/I/J), LevePriceNQ / LevePriceHQ (K/L) and the recomputed
# LeveProfitNQ / LeveProfitHQ (M/N) per affected row.

$wb = $excel.ActiveWorkbook

# ---- ALC ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H18").Value = 808
$ws.Range("I18").Value = 688
$ws.Range("J18").Value = 838
$ws.Range("K18").Value = 688
$ws.Range("L18").Value = 838
$ws.Range("M18").Value = -404
$ws.Range("N18").Value = -1406

$ws.Range("H33").Value = 163.14285
$ws.Range("J33").Value = 439.14285
$ws.Range("L33").Value = 439.14285
$ws.Range("N33").Value = -897.14285

$ws.Range("H43").Value = 2105.9412
$ws.Range("I43").Value = 3012.625
$ws.Range("J43").Value = 1300
$ws.Range("K43").Value = 3012.625
$ws.Range("L43").Value = 1300
$ws.Range("M43").Value = -2943.625
$ws.Range("N43").Value = -1438

$ws.Range("H98").Value = 1789.7368
$ws.Range("I98").Value = 2066.6
$ws.Range("J98").Value = 751.5
$ws.Range("K98").Value = 2066.6
$ws.Range("L98").Value = 751.5
$ws.Range("M98").Value = -568.5999999999999
$ws.Range("N98").Value = -3747.5

$ws.Range("H122").Value = 1789.7368
$ws.Range("I122").Value = 2066.6
$ws.Range("J122").Value = 751.5
$ws.Range("K122").Value = 6199.799999999999
$ws.Range("L122").Value = 2254.5
$ws.Range("M122").Value = -3749.799999999999
$ws.Range("N122").Value = -7154.5

# ---- ARM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2738.17
$ws.Range("I32").Value = 2738.17
$ws.Range("K32").Value = 2738.17
$ws.Range("M32").Value = -2451.17

$ws.Range("H97").Value = 3084.0688
$ws.Range("I97").Value = 3722.8572
$ws.Range("J97").Value = 1407.25
$ws.Range("K97").Value = 3722.8572
$ws.Range("L97").Value = 1407.25
$ws.Range("M97").Value = -3226.8572
$ws.Range("N97").Value = -2399.25

$ws.Range("H102").Value = 3406.4707
$ws.Range("I102").Value = 3354.6155
$ws.Range("J102").Value = 3575
$ws.Range("K102").Value = 3354.6155
$ws.Range("L102").Value = 3575
$ws.Range("M102").Value = -1732.6155
$ws.Range("N102").Value = -6819

$ws.Range("H138").Value = 90235.8
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 90235.8
$ws.Range("K138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("M138").Value = 90235.8
$ws.Range("N138").Value = -100515.8

# ---- BSM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H99").Value = 58825716
$ws.Range("I99").Value = 71430800
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 71430800
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -71429302
$ws.Range("N99").Value = -4996

# ---- CRP ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 666.6667
$ws.Range("I7").Value = 750
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 750
$ws.Range("L7").Value = 500
$ws.Range("M7").Value = -637
$ws.Range("N7").Value = -726

$ws.Range("H22").Value = 1645.8572
$ws.Range("I22").Value = 2134.2
$ws.Range("J22").Value = 425
$ws.Range("K22").Value = 2134.2
$ws.Range("L22").Value = 425
$ws.Range("M22").Value = -1784.2
$ws.Range("N22").Value = -1125

$ws.Range("H31").Value = 51457.57
$ws.Range("I31").Value = 4316.077
$ws.Range("J31").Value = 128062.5
$ws.Range("K31").Value = 4316.077
$ws.Range("L31").Value = 128062.5
$ws.Range("M31").Value = -4021.077
$ws.Range("N31").Value = -128652.5

$ws.Range("H34").Value = 51457.57
$ws.Range("I34").Value = 4316.077
$ws.Range("J34").Value = 128062.5
$ws.Range("K34").Value = 4316.077
$ws.Range("L34").Value = 128062.5
$ws.Range("M34").Value = -4114.077
$ws.Range("N34").Value = -128466.5

$ws.Range("H59").Value = 16293.286
$ws.Range("J59").Value = 16293.286
$ws.Range("L59").Value = 16293.286
$ws.Range("N59").Value = -18583.286

# ---- CUL ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 1827.36
$ws.Range("I5").Value = 231.23077
$ws.Range("J5").Value = 3556.5
$ws.Range("K5").Value = 693.69231
$ws.Range("L5").Value = 10669.5
$ws.Range("M5").Value = -581.69231
$ws.Range("N5").Value = -10893.5

$ws.Range("H135").Value = 1827.36
$ws.Range("I135").Value = 231.23077
$ws.Range("J135").Value = 3556.5
$ws.Range("K135").Value = 2081.07693
$ws.Range("L135").Value = 32008.5
$ws.Range("M135").Value = 453.9230699999998
$ws.Range("N135").Value = -37078.5

# ---- GSM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 5601.6
$ws.Range("I70").Value = 6999.5
$ws.Range("J70").Value = 4669.6665
$ws.Range("K70").Value = 6999.5
$ws.Range("L70").Value = 4669.6665
$ws.Range("M70").Value = -6729.5
$ws.Range("N70").Value = -5209.6665

$ws.Range("H73").Value = 5601.6
$ws.Range("I73").Value = 6999.5
$ws.Range("J73").Value = 4669.6665
$ws.Range("K73").Value = 6999.5
$ws.Range("L73").Value = 4669.6665
$ws.Range("M73").Value = -6063.5
$ws.Range("N73").Value = -6541.6665

$ws.Range("H136").Value = 6254.184
$ws.Range("J136").Value = 6254.184
$ws.Range("L136").Value = 18762.552
$ws.Range("N136").Value = -23862.552

# ---- LTW ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 931.129
$ws.Range("I22").Value = 268
$ws.Range("J22").Value = 1124.5416
$ws.Range("K22").Value = 268
$ws.Range("L22").Value = 1124.5416
$ws.Range("M22").Value = 27
$ws.Range("N22").Value = -1714.5416

$ws.Range("H27").Value = 931.129
$ws.Range("I27").Value = 268
$ws.Range("J27").Value = 1124.5416
$ws.Range("K27").Value = 268
$ws.Range("L27").Value = 1124.5416
$ws.Range("M27").Value = -161
$ws.Range("N27").Value = -1338.5416

$ws.Range("H55").Value = 297.53333
$ws.Range("I55").Value = 385.57144
$ws.Range("J55").Value = 220.5
$ws.Range("K55").Value = 385.57144
$ws.Range("L55").Value = 220.5
$ws.Range("M55").Value = -212.57144
$ws.Range("N55").Value = -566.5

$ws.Range("H104").Value = 15417.5
$ws.Range("J104").Value = 15417.5
$ws.Range("L104").Value = 15417.5
$ws.Range("N104").Value = -22405.5

$ws.Range("H132").Value = 2069.9375
$ws.Range("I132").Value = 2039.4062
$ws.Range("J132").Value = 2131
$ws.Range("K132").Value = 6118.2186
$ws.Range("L132").Value = 6393
$ws.Range("M132").Value = -3588.2186
$ws.Range("N132").Value = -11453

# ---- WVR ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H136").Value = 1756.0385
$ws.Range("I136").Value = 1785.5454
$ws.Range("K136").Value = 5356.6362
$ws.Range("M136").Value = -2806.6362

Write-Host "Applied $($wb.Worksheets.Count)-sheet Bahamut_Profits price/profit refresh"
